$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 172, pushing the existing rows 172-189 down to 175-192.
$ws.Rows("172:174").Insert()

# Row 172 (new)
$ws.Cells.Item(172, 1).Value = 1
$ws.Cells.Item(172, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(172, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(172, 4).Value = 45223
$ws.Cells.Item(172, 5).Value = 15
$ws.Cells.Item(172, 6).Value = 100112008
$ws.Cells.Item(172, 7).Value = "Coliflor"
$ws.Cells.Item(172, 8).Value = "Sin especificar"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 230
$ws.Cells.Item(172, 11).Value = 700
$ws.Cells.Item(172, 12).Value = 800
$ws.Cells.Item(172, 13).Value = 765
$ws.Cells.Item(172, 14).Value = "`$/unidad"
$ws.Cells.Item(172, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(172, 16).Value = 765
$ws.Cells.Item(172, 17).Value = 1
$ws.Cells.Item(172, 18).Value = "Hortaliza"

# Row 173 (new)
$ws.Cells.Item(173, 1).Value = 1
$ws.Cells.Item(173, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(173, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(173, 4).Value = 45223
$ws.Cells.Item(173, 5).Value = 15
$ws.Cells.Item(173, 6).Value = 100112008
$ws.Cells.Item(173, 7).Value = "Coliflor"
$ws.Cells.Item(173, 8).Value = "Sin especificar"
$ws.Cells.Item(173, 9).Value = "Segunda"
$ws.Cells.Item(173, 10).Value = 410
$ws.Cells.Item(173, 11).Value = 400
$ws.Cells.Item(173, 12).Value = 500
$ws.Cells.Item(173, 13).Value = 439
$ws.Cells.Item(173, 14).Value = "`$/unidad"
$ws.Cells.Item(173, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(173, 16).Value = 439
$ws.Cells.Item(173, 17).Value = 1
$ws.Cells.Item(173, 18).Value = "Hortaliza"

# Row 174 (new)
$ws.Cells.Item(174, 1).Value = 1
$ws.Cells.Item(174, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(174, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(174, 4).Value = 45223
$ws.Cells.Item(174, 5).Value = 15
$ws.Cells.Item(174, 6).Value = 100112008
$ws.Cells.Item(174, 7).Value = "Coliflor"
$ws.Cells.Item(174, 8).Value = "Sin especificar"
$ws.Cells.Item(174, 9).Value = "Tercera"
$ws.Cells.Item(174, 10).Value = 450
$ws.Cells.Item(174, 11).Value = 300
$ws.Cells.Item(174, 12).Value = 400
$ws.Cells.Item(174, 13).Value = 356
$ws.Cells.Item(174, 14).Value = "`$/unidad"
$ws.Cells.Item(174, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(174, 16).Value = 356
$ws.Cells.Item(174, 17).Value = 1
$ws.Cells.Item(174, 18).Value = "Hortaliza"

Write-Host ("Dimension now: " + $ws.UsedRange.Address())
